# Mariana_Crow_MARK.xlsx - "adds zero for NA in 2014"
#
# The 2014 column (S) was missing for the early records (rows 2-127) on the
# main "Sheet1" tab.  This fills those gaps with an explicit 0 (not seen /
# not applicable), matches what the 2015 batch of records already had, and
# refreshes the view position that was left over from the last save.  The
# backup copy of the sheet ("Sheet1 (2)") keeps a running SUM() total row
# at the bottom; re-entering that formula across the row turns it into an
# Excel "shared" formula and keeps the totals in sync.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheet1: fill column S (year 2014) with 0 for rows 2 through 127,
#    which previously had no entry at all for that year.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

$fillRange = $ws1.Range("S2:S127")
$fillRange.Value = 0

# ---------------------------------------------------------------------
# 2. Sheet1: restore the frozen header row (split after row 1) and move
#    the visible/selected cell back to L1.
# ---------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollColumn = 4
$ws1.Range("L1").Select()

# ---------------------------------------------------------------------
# 3. "Sheet1 (2)": re-apply the column totals formula across C158:S158
#    in one shot so Excel stores it as a single shared formula (it now
#    also covers the new S column).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet1 (2)")
$totalsRange = $ws4.Range("C158:S158")
$totalsRange.Formula = "=SUM(C2:C157)"

$wb.Save()
